# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect the latest scrape snapshot.
#   F2: 325 -> 327
#   F4: 1319 -> 1323
#   F5 (展览) / F6 (全部类型): 643 -> 646

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 327
    $ws.Range("F4").Value = 1323

    if ($name -eq "展览") {
        $ws.Range("F5").Value = 646
    } else {
        $ws.Range("F6").Value = 646
    }
}

$wb.Save()
